# Update job_applications sheet: replace row 2 content and append rows 3-10
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell_A2 = @'
2024-03-27 18:11:55
'@
$ws.Range("A2").Value = $cell_A2
$cell_B2 = @'
https://wellfound.com/jobs/1088986-backend-developer
'@
$ws.Range("B2").Value = $cell_B2
$cell_C2 = @'
Technical Skills:
Required strong development skills in Application Technologies listed below:
Hands on to Node.js/Nest.js technology stack.
Experience in NOSQL (Mongo) and RDBMS
Node JS, JavaScript-OOJS, ES6
Familiarity with API creation and RESTful services and Micro services
Design Pattern and its usage
Strong knowledge on scalable distributed caching technique and implementation
Strong knowledge and experience on Data Structure, Algorithm, Multi-threading and its usage
Exposure to Full stack development, deployment, feature rollout
Unit testing Framework, TDD, BDD and its implementation
Good exposure on performance optimization
Strong problem solving, debugging skills on Dev and Prod environment
Good to have:
MERN Stack with production enablement
Build tools like Grunt, Gulp, Webpack
Mocha, Chai, Jasmine, Junit/TestNG.
JWT, OAuth and Open ID Connect based Authentication
JS framework (such as React / Angular / Vue JS)
Cloud platform, preferably AWS
Bonus:
Experience with Mobile technologies
Experience with building financial applications
Exposure to Open Banking / PSD2 standard / Account Aggregator
'@
$ws.Range("C2").Value = $cell_C2
$cell_D2 = @'
Dear [Hiring Manager Name],
I'm excited to apply for the [Job Title] position at [Company Name]. As a talented full-stack developer with a passion for building innovative solutions, I’m confident that I can bring value to your team.
My strong technical skills in Node.js, MongoDB, and React, along with my expertise in cloud platforms and micro-services, align perfectly with your requirements. My commitment to performance optimization, debugging, and unit testing ensures the delivery of high-quality applications.
Past projects like Fresources and Synapse demonstrate my ability to leverage cutting-edge technologies to create user-centric solutions. Check out Synapse at [Synapse link] for a closer look at my collaborative web-app.
I’m eager to join your team and contribute to your company's success.
Thanks,
Chaitanya Anand
'@
$ws.Range("D2").Value = $cell_D2

$cell_A3 = @'
2024-03-27 18:12:10
'@
$ws.Range("A3").Value = $cell_A3
$cell_B3 = @'
https://wellfound.com/jobs/2932512-frontend-development-intern
'@
$ws.Range("B3").Value = $cell_B3
$cell_C3 = @'
Note: Just add a link to your project github/live website on the top of your application message. Let the project speak for you.
We are currently on the lookout for a dynamic and innovative front-end React developer to join our team for a 3-month internship, with the potential for extension to full-time employment based on performance. Our company operates in a fast-paced startup environment that not only demands agility and efficiency but also encourages creativity and innovation. This role is an excellent opportunity for those looking to make a tangible impact within a company and to contribute their unique ideas and solutions.
Role Overview:
As a React developer intern in our team, you will be tasked with developing and implementing responsive user interface components using React concepts.
Key Responsibilities:
Develop new user-facing features using React.js
Build reusable components and front-end libraries for future use
Optimize components for maximum performance across a vast array of web-capable devices and browsers
Collaborate with other team members
Engage in code reviews and help maintain code quality and organization
Requirements:
Proficient understanding CSS3 *Strong proficiency in JavaScript, including DOM manipulation and the JavaScript object model Deep knowledge of React.js and its core principles, especially hooks
Familiarity with RESTful APIs
Experience with common front-end development tools such as Babel, Webpack, NPM, etc.
What We Offer:
Flexible working hours with no strict timings, provided work is completed.
A stimulating work environment where you are encouraged to present your ideas and innovations
A chance to work with a team of talented individuals who are as passionate about technology as you are
Opportunity for professional growth and development, including the possibility of transitioning to full-time employment
'@
$ws.Range("C3").Value = $cell_C3
$cell_D3 = @'
I am writing to apply for the full-time opportunity at your company. I am interested in the internship but more keen on working full-time as your React developer. I have experience working on web development and front-end. I have knowledge in JavaScript and React.js.
I have built a few projects using these technologies, and I am confident that I can use my skills to contribute to your team.
[Synapse Project Link] is a productivity-enhancing collaborative web-app with combined features of Notion and Google Docs.
[Fresources Project Link] empowers college students with a platform for efficient access to curated, syllabus-relevant academic resources.
I am a hard worker and I am eager to learn new things. I am confident that I can be a valuable asset to your team.
Thank you for your time and consideration.
Yours faithfully,
Chaitanya Anand
'@
$ws.Range("D3").Value = $cell_D3

$cell_A4 = @'
2024-03-27 18:12:24
'@
$ws.Range("A4").Value = $cell_A4
$cell_B4 = @'
https://wellfound.com/jobs/2964415-frontend-developer-internship
'@
$ws.Range("B4").Value = $cell_B4
$cell_C4 = @'
Responsibilities:
Contribute to the development of new user-facing features in React.js.
Assist in designing and implementing responsive web application interfaces.
Work on building reusable code and components for future use.
Translate designs and wireframes into high-quality code.
Optimize applications for maximum speed and scalability.
Engage in learning and improving user interaction on web pages.
Collaborate with team members on various aspects of development.
Technical Skills Required:
Profound expertise in JavaScript, with a strong grasp of DOM manipulation and the JavaScript object model.
Comprehensive understanding of React.js and its fundamental principles.
Experience with popular React.js workflows, such as Flux or Redux.
Familiarity with recent EcmaScript standards.
Experience with data structure libraries (e.g., Immutable.js) is advantageous.
Knowledge of isomorphic React will be considered a plus.
Proficient in working with RESTful APIs and integrating them seamlessly into applications.
Good understanding of HTML/CSS.
Experience with GIT will be beneficial.
Awareness of modern authorization mechanisms, such as JSON Web Token.
Acquainted with modern front-end build pipelines and tools.
Skilled in using common front-end development tools like Babel, Webpack, NPM, etc.
Ability to translate business requirements into technical specifications effectively.
A talent for benchmarking and optimization.
Experience with code versioning tools (Git, SVN, Mercurial).
'@
$ws.Range("C4").Value = $cell_C4
$cell_D4 = @'
I am delighted to express my interest in the UI Developer position at your esteemed organization. As a driven and passionate student with a strong foundation in React.js and web development, I believe I possess the skills and enthusiasm to excel in this role.
My experiences at Carboledger and through my personal projects, including Fresources and Synapse [project links], have honed my abilities in designing and developing innovative user interfaces, optimizing for speed and scalability, and seamlessly integrating APIs. I am eager to contribute my expertise in React.js to your team and drive the creation of exceptional user experiences.
'@
$ws.Range("D4").Value = $cell_D4

$cell_A5 = @'
2024-03-27 18:12:37
'@
$ws.Range("A5").Value = $cell_A5
$cell_B5 = @'
https://wellfound.com/jobs/2964456-intern-software-development
'@
$ws.Range("B5").Value = $cell_B5
$cell_C5 = @'
We are seeking a highly motivated and talented intern to join our team as a React Native / GraphQL / React software engineering intern. In this role, you will have the opportunity to work on a product to develop and host a GraphQL API for managing user data using React Native and React. This internship will provide valuable hands-on experience in React Native, GraphQL development, API design, and component-based engineering.
Selected intern's day-to-day responsibilities include:
Collaborate with the development team to implement the GraphQL schema and resolvers using React Native/React
Assist in testing and debugging the GraphQL API, React Native components, and front-end Application to ensure functionality and performance
Contribute to the documentation of the API endpoints, schema, usage examples, and React components
Work closely with senior developers to learn and grow your skills in React Native, GraphQL, and React development
Requirements:
Currently pursuing a Bachelor's or Master's degree in Computer Science, Software Engineering, or a related field
Familiarity with UI development concepts and technologies, such as React.js, React Native, and GraphQL
Basic knowledge of GraphQL, React Native, and React and their concepts is a plus, but not required
Excellent problem-solving skills and attention to detail
Ability to work independently as well as collaborate effectively in a team environment
Strong communication skills and willingness to learn new technologies
Benefits:
Opportunity to work closely with experienced developers and learn best practices in software engineering
Flexible work hours and remote work options
Potential for future employment opportunities based on performance and business needs
Duration:
This internship position is available for 6 months with the possibility of extension based on performance and mutual interest.
'@
$ws.Range("C5").Value = $cell_C5
$cell_D5 = @'
As a college student brimming with enthusiasm and technical prowess, I'm eager to join your team as an intern to delve into React Native, GraphQL, and React development. I'm particularly intrigued by the opportunity to work on your GraphQL API and leverage my skills in React Native/React.
My projects, including Fresources [https://fresources.tech] and Synapse [https://github.com/UltigendLemate/synapse], showcase my ability to design, develop, and collaborate effectively in a team environment. I'm confident that my dedication and passion for software engineering will make me a valuable asset to your organization.
Thank you for considering my application. I look forward to hearing from you soon.
'@
$ws.Range("D5").Value = $cell_D5

$cell_A6 = @'
2024-03-27 18:12:50
'@
$ws.Range("A6").Value = $cell_A6
$cell_B6 = @'
https://wellfound.com/jobs/2830594-full-stack-developer
'@
$ws.Range("B6").Value = $cell_B6
$cell_C6 = @'
Please read carefully before applying:
This role is remote and can be freelance. Open to full-time too. There is a lot of flexibility for motivated talent.
Please submit github and beehance or similar repos to demonstrate frontend and development skills.
Skills:
Strong frontend skills - JS, frameworks, HTML, CSS
Experience in Python, serious programming, not some Jupyter notebook project
Rest during interview process
'@
$ws.Range("C6").Value = $cell_C6
$cell_D6 = @'
I'm a final-year student who lives and breathes frontend. I've built a platform used by over 30K students, have built collaborative tools, and even turned the power of AI into a market research assistant.
I've got a knack for JS, frameworks, HTML, and CSS. And when it comes to Python, I'm not just a Jupyter notebook warrior.
You can find all my work and links here:
- Website: [fresources.tech](https://fresources.tech)
- Github: [github.com/UltigendLemate](https://github.com/UltigendLemate)
- LinkedIn: [linkedin.com/in/chaitanyanand](https://linkedin.com/in/chaitanyanand)
'@
$ws.Range("D6").Value = $cell_D6

$cell_A7 = @'
2024-03-27 18:13:03
'@
$ws.Range("A7").Value = $cell_A7
$cell_B7 = @'
https://wellfound.com/jobs/2627836-full-stack-developer
'@
$ws.Range("B7").Value = $cell_B7
$cell_C7 = @'
Job Description:
We are seeking a highly skilled and experienced Backend Developer with expertise in Python and Nodejs Development to join our team. As a Backend Developer, you will be responsible for designing and developing robust, scalable, and secure backend systems for our applications.
Requirements:
Architect, develop, and maintain a highly scalable backend platform using Node.js and Python.
Collaborate with frontend developers to ensure seamless integration and responsiveness across the platform.
Craft efficient, maintainable, and scalable code, adhering to best practices and coding standards.
Continuously optimize platform performance and scalability, identifying and addressing bottlenecks and inefficiencies.
Tackle complex technical issues, including troubleshooting and debugging.
Work closely with DevOps and infrastructure teams to ensure the platform's reliability, security, and scalability.
Lead code reviews and mentor junior team members, fostering a culture of excellence.
Stay at the forefront of industry trends and emerging technologies, suggesting innovative solutions to enhance our platform's scalability.
Qualifications:
Bachelor's degree in Computer Science, Software Engineering, or a related field (or equivalent work experience).
Proven track record as a backend developer, with a strong focus on Node.js and Python.
Deep expertise in server-side development, asynchronous programming, and RESTful APIs.
Proficiency in database systems (e.g., SQL, NoSQL) and advanced data modeling.
Experience with cloud platforms (e.g., AWS, Azure, GCP) and containerization (e.g., Docker, Kubernetes) is highly desirable.
Exceptional problem-solving skills and the ability to troubleshoot intricate issues.
Strong communication and teamwork skills.
A passion for learning and staying abreast of industry trends and technologies.
'@
$ws.Range("C7").Value = $cell_C7
$cell_D7 = @'
With deep expertise in Python and Node.js, my backend development skills stand out. My projects Fresources, Synapse, and Quik Planr demonstrate my ability to design scalable platforms and optimize performance. Eager to contribute to your team's success.
[fresources.tech, github.com/UltigendLemate/synapse, github.com/UltigendLemate/code-kshetra]
'@
$ws.Range("D7").Value = $cell_D7

$cell_A8 = @'
2024-03-27 18:13:16
'@
$ws.Range("A8").Value = $cell_A8
$cell_B8 = @'
https://wellfound.com/jobs/2925735-frontend-developer
'@
$ws.Range("B8").Value = $cell_B8
$cell_C8 = @'
Job Overview:
We are looking for a talented and experienced Frontend Developer to join our team and contribute to the development of our project. The ideal candidate will have a strong background in frontend development, proficiency in modern web technologies, and a keen eye for design. As a Frontend Developer, you will work closely with our backend developers, UI/UX designers, and project managers to deliver intuitive and responsive user interfaces for our web application.
Responsibilities:
Develop user-facing features and interfaces for our web application using HTML, CSS, and JavaScript.
Collaborate with UI/UX designers to translate design mockups and wireframes into high-quality frontend code.
Ensure the technical feasibility of UI/UX designs and implement responsive design principles to deliver a seamless user experience across devices.
Integrate frontend components with backend APIs and services to support dynamic data rendering and user interactions.
Optimize application performance and load times by implementing best practices and performance optimization techniques.
Conduct code reviews, debugging sessions, and testing to identify and resolve frontend issues and ensure code quality and stability.
Stay updated with emerging frontend technologies, frameworks, and trends to continuously improve our frontend architecture and development processes.
Communicate effectively with team members, project managers, and stakeholders to understand project requirements and provide timely updates on progress and challenges.
Requirements:
Proven experience as a Frontend Developer or similar role, with a portfolio showcasing frontend projects and applications.
Proficiency in frontend technologies including HTML5, CSS3, JavaScript (ES6+), and modern JavaScript frameworks/libraries such as React, Vue.js, or Angular.
Strong understanding of responsive web design principles, cross-browser compatibility, and accessibility standards.
Experience with frontend build tools and package managers such as Webpack, Babel, npm, or Yarn.
Familiarity with version control systems (e.g., Git) and collaborative development workflows.
Knowledge of UI/UX design principles and experience working with design tools such as Adobe XD, Sketch, or Figma.
Excellent problem-solving skills and attention to detail, with the ability to debug frontend issues and optimize performance.
Strong communication and collaboration skills, with the ability to work effectively in a team environment.
Experience with backend technologies (e.g., Node.js, Python, Java) and RESTful APIs is a plus but not required.
Availability to work with our team on a full-time basis, with flexible hours and the ability to meet project deadlines.
'@
$ws.Range("C8").Value = $cell_C8
$cell_D8 = @'
I'm Chaitanya Anand, a full-stack developer with a track record of building innovative web applications. I'm eager to join your frontend team and contribute to your project.
I have a strong portfolio of frontend projects, including Fresources (fresources.tech), Synapse (github.com/UltigendLemate/synapse), and Quik Planr (github.com/UltigendLemate/code-kshetra).
I'm proficient in modern frontend technologies and have experience working with UI/UX designers and backend developers.
I'm confident that I can make a valuable contribution to your team and help you deliver high-quality frontend solutions.
'@
$ws.Range("D8").Value = $cell_D8

$cell_A9 = @'
2024-03-27 18:13:29
'@
$ws.Range("A9").Value = $cell_A9
$cell_B9 = @'
https://wellfound.com/jobs/2925728-backend-developers
'@
$ws.Range("B9").Value = $cell_B9
$cell_C9 = @'
We are seeking a skilled Backend Developer to join our team and help build the backend infrastructure for our exciting project. The ideal candidate will have a strong background in backend development, experience with cloud deployment, Google Maps API integration, authentication mechanisms, and a passion for creating scalable and efficient systems. As a Backend Developer, you will collaborate with frontend developers, UI/UX designers, and other team members to deliver high-quality software solutions.
Responsibilities:
Design, develop, and maintain the backend infrastructure for our project.
Implement scalable and efficient server-side logic using appropriate programming languages and frameworks.
Integrate with databases, APIs, and external services to support frontend functionalities.
Optimize application performance and scalability by implementing best practices and optimization techniques.
Collaborate with frontend developers to define API specifications and ensure seamless integration between frontend and backend components.
Implement cloud deployment strategies and manage infrastructure on platforms like AWS, Azure, or Google Cloud.
Integrate Google Maps API for location-based services and route optimization.
Implement authentication and authorization mechanisms to ensure secure access to our application.
Participate in code reviews, debugging sessions, and troubleshooting efforts to identify and resolve issues.
Work closely with UI/UX designers and project managers to understand project requirements and deliver solutions that meet business objectives.
Stay updated with emerging technologies, trends, and best practices in backend development to continuously improve our systems.
Requirements:
Bachelor's degree in Computer Science, Engineering, or a related field (or equivalent work experience).
Proven experience as a Backend Developer or similar role, with a track record of delivering backend solutions for web applications.
Proficiency in backend programming languages such as Python, Node.js, or Java.
Strong understanding of web development concepts, including RESTful APIs, MVC architecture, and server-side templating languages.
Experience with database management systems (e.g., MySQL, PostgreSQL, MongoDB) and proficiency in writing complex SQL queries.
Familiarity with cloud platforms (e.g., AWS, Azure, Google Cloud) and experience with infrastructure as code (e.g., Terraform, CloudFormation).
Experience integrating Google Maps API for location-based services, route optimization, and geocoding.
Familiarity with authentication and authorization mechanisms such as OAuth 2.0, JWT, or session-based authentication.
Knowledge of version control systems (e.g., Git) and collaborative development workflows.
Excellent problem-solving skills and attention to detail, with the ability to analyze and debug complex issues.
Strong communication and collaboration skills, with the ability to work effectively in a team environment.
Knowledge of frontend technologies (e.g., HTML, CSS, JavaScript) and experience with frontend frameworks (e.g., React, Angular, Vue.js) is desirable but not required.
'@
$ws.Range("C9").Value = $cell_C9
$cell_D9 = @'
My projects Synapse and Fresources showcase my proficiency in backend development and API integration. I am eager to apply my skills to your project, leveraging them to build a scalable and secure backend infrastructure. I am a quick learner and an effective communicator, and I am confident in my ability to make a significant contribution to your team.
Here are the links to my projects:
- Synapse: [LINK]
- Fresources: [LINK]
'@
$ws.Range("D9").Value = $cell_D9

$cell_A10 = @'
2024-03-27 18:13:43
'@
$ws.Range("A10").Value = $cell_A10
$cell_B10 = @'
https://wellfound.com/jobs/2866106-software-engineering-intern-web-dev
'@
$ws.Range("B10").Value = $cell_B10
$cell_C10 = @'
Role Description
As an SDE Intern at BigCircle, you will play a critical role in the development of innovative software solutions and SaaS products. You will be responsible for designing, coding, and implementing both the front-end and back-end of applications, ensuring seamless integration and exceptional user experience. Your expertise in multiple programming languages and platforms will contribute to our mission of revolutionizing the software development industry in India.
Responsibilities:
Develop, test, and maintain both front-end and back-end components of software applications and SaaS products.
Collaborate with cross-functional teams, including designers, product managers, and other developers, to deliver high-quality solutions according to client requirements.
Write clean, efficient, and maintainable code, ensuring adherence to coding standards and best practices.
Participate in code reviews to identify and resolve any issues or bugs.
Optimize application performance through continuous monitoring, troubleshooting, and applying performance improvement techniques.
Stay updated with the latest industry trends and technologies, and proactively suggest innovative solutions to enhance the development process.
Collaborate with the Quality Assurance team to ensure the delivery of bug-free and high-quality software.
Contribute to the growth and development of the software ecosystem in India by actively participating in knowledge-sharing sessions
Qualifications
Strong understanding of Computer Science fundamentals
Programming proficiency in relevant languages such as JavaScript/TypeScript
Experience with React/Node.js is a plus.
Understanding of Object-Oriented Programming (OOP)
Experience using development tools such as Git, Visual Studio Code
Problem-solving skills and ability to troubleshoot code effectively
Excellent written and verbal communication skills
Pursuing or recently completed a Bachelor's or Master's degree in Computer Science or related field
'@
$ws.Range("C10").Value = $cell_C10
$cell_D10 = @'
Dear [Hiring Manager name],
As a budding developer eager to revolutionize the software industry, I'm thrilled to apply as an SDE Intern at BigCircle. My hands-on experience in developing end-to-end solutions, backed by my academic grounding in Computer Science, makes me an ideal fit for this opportunity.
My projects like Fresources and Synapse showcase my ability to build user-friendly, collaborative web applications from scratch. I am also adept at using React/Node.js and have a solid grasp of OOP principles.
I am confident that I can leverage my skills to contribute to BigCircle's cutting-edge solutions. I am eager to apply my knowledge to real-world problems and grow as a developer in your esteemed organization.
Let's explore how I can revolutionize the software industry alongside BigCircle. I am excited to learn more about this internship and discuss how my expertise can benefit your team.
Thank you for your time and consideration.
Best regards,
Chaitanya Anand
'@
$ws.Range("D10").Value = $cell_D10

